# Correct the "breakdown comment" text for the volume-type metrics: these
# cells previously said things like "number of businesses" / "employment
# volume", which read oddly in the sentence template used for the chart's
# breakdown caption. They are reworded as "share of ..." / "... share" so
# the generated sentence reads correctly.
#
# Setting the cell .Value directly is sufficient: each of the old strings
# being replaced is referenced by exactly one cell in the workbook, so
# Excel drops the now-unused shared-string entries and appends the new
# text at the end of the shared-string table on save - matching the
# target edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L11").Value = "share of businesses"
$ws.Range("L15").Value = "FE participation share"
$ws.Range("L14").Value = "share of FE achievements"
$ws.Range("L10").Value = "share of online job adverts"
$ws.Range("L6").Value = "employment volume share"

# Restore the view to what was saved with the workbook: scrolled back to
# the top (no frozen/offset topLeftCell) with a single cell (F3) selected,
# rather than the previous B2:B9 multi-cell selection scrolled to row 17.
[void]$ws.Range("F3").Select()
